$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "1.0210 at -2.49"
$ws.Range("C2").Value = "1.0420 at -121.72"
$ws.Range("D2").Value = "1.0174 at 117.83"

$ws.Range("C3").Value = "1.0328 at -121.90"
$ws.Range("D3").Value = "1.0154 at 117.86"

$ws.Range("C4").Value = "1.0311 at -121.98"
$ws.Range("D4").Value = "1.0134 at 117.90"

$ws.Range("B5").Value = "1.0180 at -2.55"
$ws.Range("C5").Value = "1.0401 at -121.77"
$ws.Range("D5").Value = "1.0148 at 117.83"

$ws.Range("B6").Value = "0.9940 at -3.23"
$ws.Range("C6").Value = "1.0218 at -122.22"
$ws.Range("D6").Value = "0.9960 at 117.35"

$ws.Range("B7").Value = "0.9900 at -5.30"
$ws.Range("C7").Value = "1.0529 at -122.34"
$ws.Range("D7").Value = "0.9777 at 116.03"

$ws.Range("B8").Value = "0.9881 at -5.32"
$ws.Range("D8").Value = "0.9757 at 115.93"

$ws.Range("D9").Value = "0.9737 at 115.78"

$ws.Range("B10").Value = "0.9835 at -5.55"
$ws.Range("C10").Value = "1.0553 at -122.52"
$ws.Range("D10").Value = "0.9758 at 116.04"

$ws.Range("B11").Value = "0.9900 at -5.30"
$ws.Range("C11").Value = "1.0529 at -122.34"
$ws.Range("D11").Value = "0.9777 at 116.03"

$ws.Range("B12").Value = "0.9825 at -5.25"
